# Auto-generated edit script applying numeric corrections to Ixion Profits workbook
# (scheduled price-refresh run: recompute H-N market/profit columns across sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 883.9483
$ws.Range("J129").Value = 923.4259
$ws.Range("L129").Value = 2770.2777
$ws.Range("N129").Value = -12770.2777
$ws.Range("H137").Value = 1590
$ws.Range("I137").Value = 1030.4706
$ws.Range("J137").Value = 1929.7142
$ws.Range("K137").Value = 3091.4118
$ws.Range("L137").Value = 5789.142599999999
$ws.Range("M137").Value = -541.4118000000003
$ws.Range("N137").Value = -10889.1426
$ws.Range("H138").Value = 2538.9775
$ws.Range("I138").Value = 997.4524
$ws.Range("J138").Value = 3916.5107
$ws.Range("K138").Value = 2992.3572
$ws.Range("L138").Value = 11749.5321
$ws.Range("M138").Value = 2147.6428
$ws.Range("N138").Value = -22029.5321

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2555
$ws.Range("I32").Value = 1694.119
$ws.Range("K32").Value = 1694.119
$ws.Range("M32").Value = -1407.119
$ws.Range("H74").Value = 1210.2452
$ws.Range("I74").Value = 953.72974
$ws.Range("J74").Value = 1803.4375
$ws.Range("K74").Value = 953.72974
$ws.Range("L74").Value = 1803.4375
$ws.Range("M74").Value = -79.72973999999999
$ws.Range("N74").Value = -3551.4375
$ws.Range("H77").Value = 1210.2452
$ws.Range("I77").Value = 953.72974
$ws.Range("J77").Value = 1803.4375
$ws.Range("K77").Value = 4768.6487
$ws.Range("L77").Value = 9017.1875
$ws.Range("M77").Value = -400.6486999999997
$ws.Range("N77").Value = -17753.1875
$ws.Range("H132").Value = 4904.1934
$ws.Range("I132").Value = 2347.8572
$ws.Range("J132").Value = 7009.4116
$ws.Range("K132").Value = 7043.571599999999
$ws.Range("L132").Value = 21028.2348
$ws.Range("M132").Value = -4513.571599999999
$ws.Range("N132").Value = -26088.2348

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4183.942
$ws.Range("I31").Value = 1967.65
$ws.Range("J31").Value = 5088.551
$ws.Range("K31").Value = 1967.65
$ws.Range("L31").Value = 5088.551
$ws.Range("M31").Value = -1672.65
$ws.Range("N31").Value = -5678.551
$ws.Range("H34").Value = 4183.942
$ws.Range("I34").Value = 1967.65
$ws.Range("J34").Value = 5088.551
$ws.Range("K34").Value = 1967.65
$ws.Range("L34").Value = 5088.551
$ws.Range("M34").Value = -1765.65
$ws.Range("N34").Value = -5492.551
$ws.Range("H58").Value = 1888.8235
$ws.Range("I58").Value = 1222.4445
$ws.Range("J58").Value = 2638.5
$ws.Range("K58").Value = 1222.4445
$ws.Range("L58").Value = 2638.5
$ws.Range("M58").Value = -1019.4445
$ws.Range("N58").Value = -3044.5
$ws.Range("H136").Value = 1888.8235
$ws.Range("I136").Value = 1222.4445
$ws.Range("J136").Value = 2638.5
$ws.Range("K136").Value = 3667.3335
$ws.Range("L136").Value = 7915.5
$ws.Range("M136").Value = -1117.3335
$ws.Range("N136").Value = -13015.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3212.5
$ws.Range("I68").Value = 4653.778
$ws.Range("J68").Value = 1870.6207
$ws.Range("K68").Value = 13961.334
$ws.Range("L68").Value = 5611.8621
$ws.Range("M68").Value = -13150.334
$ws.Range("N68").Value = -7233.8621
$ws.Range("H71").Value = 3212.5
$ws.Range("I71").Value = 4653.778
$ws.Range("J71").Value = 1870.6207
$ws.Range("K71").Value = 41884.002
$ws.Range("L71").Value = 16835.5863
$ws.Range("M71").Value = -37828.002
$ws.Range("N71").Value = -24947.5863
$ws.Range("H130").Value = 8662.5
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 8662.5
$ws.Range("K130").Value = 0
$ws.Range("L130").ClearContents()
$ws.Range("M130").Value = 25987.5
$ws.Range("N130").Value = -36027.5
$ws.Range("H131").Value = 20371466
$ws.Range("J131").Value = 26317150
$ws.Range("L131").Value = 78951450
$ws.Range("N131").Value = -78961530
$ws.Range("H137").Value = 25654244
$ws.Range("I137").Value = 1155.4615
$ws.Range("J137").Value = 38480788
$ws.Range("K137").Value = 3466.3845
$ws.Range("L137").Value = 115442364
$ws.Range("M137").Value = 1633.6155
$ws.Range("N137").Value = -115452564

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 10000000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H70").Value = 5387.517
$ws.Range("I70").Value = 5757.6055
$ws.Range("J70").Value = 4684.35
$ws.Range("K70").Value = 5757.6055
$ws.Range("L70").Value = 4684.35
$ws.Range("M70").Value = -5487.6055
$ws.Range("N70").Value = -5224.35
$ws.Range("H73").Value = 5387.517
$ws.Range("I73").Value = 5757.6055
$ws.Range("J73").Value = 4684.35
$ws.Range("K73").Value = 5757.6055
$ws.Range("L73").Value = 4684.35
$ws.Range("M73").Value = -4821.6055
$ws.Range("N73").Value = -6556.35
$ws.Range("H122").Value = 4594295.5
$ws.Range("I122").Value = 2947880
$ws.Range("J122").Value = 16668009
$ws.Range("K122").Value = 8843640
$ws.Range("L122").Value = 50004027
$ws.Range("M122").Value = -8841190
$ws.Range("N122").Value = -50008927
$ws.Range("H126").Value = 7975.125
$ws.Range("I126").Value = 8900.857
$ws.Range("J126").Value = 1495
$ws.Range("K126").Value = 26702.571
$ws.Range("L126").Value = 4485
$ws.Range("M126").Value = -24232.571
$ws.Range("N126").Value = -9425

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 15874289
$ws.Range("I46").Value = 23810398
$ws.Range("J46").Value = 2071.4285
$ws.Range("K46").Value = 23810398
$ws.Range("L46").Value = 2071.4285
$ws.Range("M46").Value = -23810210
$ws.Range("N46").Value = -2447.4285
$ws.Range("H55").Value = 18518940
$ws.Range("I55").Value = 338.23077
$ws.Range("J55").Value = 35714784
$ws.Range("K55").Value = 338.23077
$ws.Range("L55").Value = 35714784
$ws.Range("M55").Value = -165.23077
$ws.Range("N55").Value = -35715130
$ws.Range("H132").Value = 21671162
$ws.Range("I132").Value = 28893496
$ws.Range("J132").Value = 4159.8
$ws.Range("K132").Value = 86680488
$ws.Range("L132").Value = 12479.4
$ws.Range("M132").Value = -86677958
$ws.Range("N132").Value = -17539.4
$ws.Range("H136").Value = 6429.357
$ws.Range("I136").Value = 3205.6843
$ws.Range("K136").Value = 9617.052899999999
$ws.Range("M136").Value = -7067.052899999999

Write-Host "Applied 160 cell updates across 6 sheets"
